# Apply cryptos.xlsx diff via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '58.285.24'
$ws.Range('E2').Value = '  -0.47%  '

$ws.Range('D3').Value = '3.141.78'
$ws.Range('E3').Value = '  +1.50%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.27%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.87%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '3.143.26'
$ws.Range('E8').Value = '  +1.52%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.445'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.86%  '

$ws.Range('E10').Value = '  -2.53%  '

$ws.Range('E11').Value = '  +0.42%  '

$ws.Range('E12').Value = '  +2.14%  '

$ws.Range('D13').Value = '3.682.08'
$ws.Range('E13').Value = '  +1.60%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.134'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.78%  '

$ws.Range('E16').Value = '  -0.35%  '

$ws.Range('D17').Value = '58.307.63'
$ws.Range('E17').Value = '  -0.47%  '

$ws.Range('D18').Value = '3.142.74'
$ws.Range('E18').Value = '  +1.67%  '

$ws.Range('E19').Value = '  -0.48%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.80%  '

$ws.Range('E21').Value = '  -1.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '343.59'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.81%  '

$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.513'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.169'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.67%  '

$ws.Range('E27').Value = '  -0.06%  '

$ws.Range('D28').Value = '0.0₃0933'
$ws.Range('E28').Value = '  +2.03%  '

$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.86%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.89'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.11'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.41%  '

$ws.Range('E34').Value = '  -0.57%  '

$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.81'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.22%  '

$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.43'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.68%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.21'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.23%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.92%  '

$ws.Range('E39').Value = '  -5.22%  '

$ws.Range('E40').Value = '  +12.17%  '

$ws.Range('E41').Value = '  -1.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.709'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.72%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.52%  '

$ws.Range('D44').Value = '3.179.80'
$ws.Range('E44').Value = '  +1.41%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.27%  '

$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0264'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.89%  '

$ws.Range('D48').Value = '2.276.45'
$ws.Range('E48').Value = '  -0.10%  '

$ws.Range('E49').Value = '  +5.09%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.88%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.52%  '
